$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction13")

$values = @(2,3,4,5,6,7,8,10,11,13,14,15,16,17,18,19,20,21,23,27,29,30,34,37,38,39,40,41,42,43,44,45,46,47,48,50,52,54,55,56,58,59,60,62,63,64,65,66,67,68,69,70,71,72,73,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $values[$i]
}
